$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.133.25'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '3.123.27'
$ws.Range("E3").Value = '  +0.77%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("D6").Value = "'" + '174.60'
$ws.Range("E6").Value = '  +1.10%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -0.27%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").Value = "'" + '6.41'
$ws.Range("E10").Value = '  -0.65%  '
$ws.Range("D11").Value = "'" + '0.480'
$ws.Range("E11").Value = '  -0.73%  '
$ws.Range("D12").Value = "'" + '0.0000250'
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("D13").Value = "'" + '37.24'
$ws.Range("E13").Value = '  -0.57%  '
$ws.Range("E14").Value = '  -1.56%  '
$ws.Range("D15").Value = '3.641.45'
$ws.Range("E15").Value = '  +0.74%  '
$ws.Range("D16").Value = '67.124.42'
$ws.Range("E16").Value = '  +0.38%  '
$ws.Range("D17").Value = "'" + '7.12'
$ws.Range("E17").Value = '  -0.80%  '
$ws.Range("D18").Value = '3.124.73'
$ws.Range("E18").Value = '  +0.72%  '
$ws.Range("D19").Value = "'" + '16.44'
$ws.Range("E19").Value = '  +1.25%  '
$ws.Range("D20").Value = "'" + '492.62'
$ws.Range("E20").Value = '  +2.12%  '
$ws.Range("D21").Value = "'" + '7.96'
$ws.Range("E21").Value = '  +5.87%  '
$ws.Range("E22").Value = '  -0.90%  '
$ws.Range("D23").Value = "'" + '84.21'
$ws.Range("E23").Value = '  +0.25%  '
$ws.Range("D24").Value = "'" + '13.20'
$ws.Range("E24").Value = '  +0.62%  '
$ws.Range("E25").Value = '  -2.74%  '
$ws.Range("D26").Value = "'" + '10.40'
$ws.Range("E26").Value = '  +3.50%  '
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("E28").Value = '  -0.86%  '
$ws.Range("D29").Value = "'" + '2.36'
$ws.Range("E29").Value = '  -1.79%  '
$ws.Range("E30").Value = '  -0.29%  '
$ws.Range("D31").Value = "'" + '28.64'
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("E32").Value = '  -0.34%  '
$ws.Range("D33").Value = '0.0₃0952'
$ws.Range("E33").Value = '  -5.59%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("D35").Value = "'" + '5.88'
$ws.Range("E35").Value = '  -0.40%  '
$ws.Range("E36").Value = '  -1.49%  '
$ws.Range("D37").Value = "'" + '47.17'
$ws.Range("E37").Value = '  -1.82%  '
$ws.Range("D38").Value = "'" + '2.06'
$ws.Range("E38").Value = '  -2.56%  '
$ws.Range("D39").Value = "'" + '0.311'
$ws.Range("E39").Value = '  -1.49%  '
$ws.Range("E40").Value = '  +1.79%  '
$ws.Range("D41").Value = "'" + '8.52'
$ws.Range("E41").Value = '  -1.69%  '
$ws.Range("D42").Value = '2.820.23'
$ws.Range("E42").Value = '  -0.58%  '
$ws.Range("D43").Value = "'" + '2.61'
$ws.Range("E43").Value = '  -7.24%  '
$ws.Range("D44").Value = "'" + '383.79'
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("E45").Value = '  -2.14%  '
$ws.Range("D46").Value = "'" + '135.49'
$ws.Range("E46").Value = '  +0.65%  '
$ws.Range("D48").Value = "'" + '24.95'
$ws.Range("E48").Value = '  -0.04%  '
$ws.Range("E49").Value = '  -0.85%  '
$ws.Range("E50").Value = '  -0.63%  '
$ws.Range("D51").Value = "'" + '6.75'
$ws.Range("E51").Value = '  -1.41%  '
